$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "vinhtranak02092kz@gmail.com"
$ws.Range("B2").Value = "Raul123"
$ws.Range("A3").Value = "vinhtranak02092k@gmail.com"
$ws.Range("B3").Value = "Raul1231"
$ws.Range("A4").Value = "vinhtranak02092kz@gmail.com"
$ws.Range("B4").Value = "Raul1231"

$ws.Range("E41").Select()
